$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds locale-formatted text (grouping dots, fixed
# trailing zeros), not real numbers, in the source sheet. Every new Price
# value below is written with a leading apostrophe so Excel keeps storing
# it as text instead of silently reparsing it as a number.

$ws.Range("D2").Value = "'36.881.02"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "'2.033.51"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'244.29"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").Value = "'0.653"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("D7").Value = "'57.58"
$ws.Range("E7").Value = "  -2.71%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").Value = "'0.0765"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("D12").Value = "'15.20"
$ws.Range("E12").Value = "  -5.76%  "
$ws.Range("D13").Value = "'0.874"
$ws.Range("E13").Value = "  +7.14%  "
$ws.Range("D14").Value = "'2.329.87"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "'2.041.86"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "'18.05"
$ws.Range("E17").Value = "  +4.40%  "
$ws.Range("D18").Value = "'36.809.70"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "'73.20"
$ws.Range("E19").Value = "  -2.15%  "
$ws.Range("D20").Value = "'0.0₃0879"
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "'234.61"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").Value = "'9.50"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").Value = "'167.52"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "'2.12"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").Value = "'19.78"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").Value = "'5.47"
$ws.Range("E29").Value = "  +15.49%  "
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("E31").Value = "  -5.00%  "
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").Value = "'0.0607"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'0.0862"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").Value = "'1.84"
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -5.46%  "
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("D40").Value = "'5.16"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").Value = "'0.0220"
$ws.Range("E41").Value = "  -0.51%  "

# Rows 42/43 swap places: ARBITRUM moves up to row 42, Cronos moves down to
# row 43, each carrying its own refreshed Price/Volume(1h) figures.
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.13"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0943"
$ws.Range("E43").Value = "  -13.81%  "

$ws.Range("D44").Value = "'96.44"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'16.74"
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("D46").Value = "'1.285.40"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'2.33"
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "'3.62"
$ws.Range("E49").Value = "  +3.24%  "

# Rows 50/51 swap places: RocketPoolETH moves up to row 50, FraxShare moves
# down to row 51, each carrying its own refreshed Price/Volume(1h) figures.
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "'2.215.47"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'6.64"
$ws.Range("E51").Value = "  -2.35%  "
